# Updated cryptos list on Sat Jun 22 19:22:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.287.12"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.496.01"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "586.24"
$ws.Range("E5").Value = "  +0.23%  "

# Row 6 - Solana
$ws.Range("D6").Value = "134.29"
$ws.Range("E6").Value = "  +1.68%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.488"
$ws.Range("E8").Value = "  +0.69%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.08%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +2.13%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +2.19%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.087.35"
$ws.Range("E12").Value = "  -0.58%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +1.76%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.04%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.492.94"
$ws.Range("E15").Value = "  -0.41%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "25.80"
$ws.Range("E16").Value = "  -6.36%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "64.348.45"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "9.89"
$ws.Range("E18").Value = "  +0.72%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "5.75"
$ws.Range("E19").Value = "  +2.36%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  -3.45%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "395.54"
$ws.Range("E21").Value = "  +2.96%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "0.570"
$ws.Range("E22").Value = "  -0.81%  "

# Row 23 - WrappedeETH
$ws.Range("D23").Value = "3.634.03"
$ws.Range("E23").Value = "  -0.66%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "74.41"
$ws.Range("E24").Value = "  +0.92%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  -0.02%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "7.40"
$ws.Range("E28").Value = "  -1.50%  "

# Rows 29-31 rotate: old31 -> new29, old29 -> new30, old30 -> new31
# Row 29 - InternetComputer(DFINITY)
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "8.27"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30 - Fetch.AI
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "1.49"
$ws.Range("E30").Value = "  -5.25%  "

# Row 31 - PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  -0.18%  "

# Row 32 - RenzoRestakedETH
$ws.Range("D32").Value = "3.514.46"
$ws.Range("E32").Value = "  -0.41%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  +3.47%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  +0.04%  "

# Row 35 - EthereumClassic
$ws.Range("D35").Value = "23.42"
$ws.Range("E35").Value = "  -0.37%  "

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").Value = "  -3.47%  "

# Row 37 - Aptos
$ws.Range("D37").Value = "6.90"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -0.74%  "

# Row 39 - Monero
$ws.Range("D39").Value = "166.07"
$ws.Range("E39").Value = "  +4.42%  "

# Row 40 - Hedera
$ws.Range("D40").Value = "0.0783"
$ws.Range("E40").Value = "  -0.81%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  -0.88%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.01%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "25.23"
$ws.Range("E43").Value = "  -4.09%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "4.40"
$ws.Range("E44").Value = "  -0.29%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +2.40%  "

# Row 46 - ONDO
$ws.Range("E46").Value = "  -3.55%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.460.18"
$ws.Range("E47").Value = "  +0.68%  "

# Row 48 - Cosmos
$ws.Range("E48").Value = "  -0.85%  "

# Row 49 - SuiNetwork
$ws.Range("D49").Value = "0.899"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -0.86%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "21.18"
$ws.Range("E51").Value = "  -1.32%  "
